$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width (Name column) ---
$ws.Columns("B").ColumnWidth = 22.3

# --- Existing UV_1 / UV_2 state rows: Source column now points to the
#     blockade word address instead of the (now unused) state word ---
$ws.Range("C35").Value = "DB10.DBW14"
$ws.Range("C37").Value = "DB11.DBW14"

# --- Rows 38-42 (AI_1..LI_2): AlarmLimitMax now enabled ---
$ws.Range("K38:K42").Value = 1

# --- New rows 43-48: valve station Mode / Open-Close / Blockade tags ---
# Pre-apply the same "vertical center" formatting used by the rest of the
# data rows so untouched cells in the new rows match the existing style.
$ws.Range("A43:L48").VerticalAlignment = -4108

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "UV_1_MODE"
$ws.Range("C43").Value = "DB10.DBX12.1"
$ws.Range("D43").Value = "BOOL"
$ws.Range("E43").Value = "None"
$ws.Range("F43").Value = "0 - Auto, 1 -  Manual"
$ws.Range("G43").Value = "False"
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1
$ws.Range("L43").Value = "False"

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "UV_2_MODE"
$ws.Range("C44").Value = "DB11.DBX12.1"
$ws.Range("D44").Value = "BOOL"
$ws.Range("E44").Value = "None"
$ws.Range("F44").Value = "0 - Auto, 1 -  Manual"
$ws.Range("G44").Value = "False"
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1
$ws.Range("L44").Value = "False"

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "UV_1_OPEN_CLOSE"
$ws.Range("C45").Value = "DB10.DBX12.0"
$ws.Range("D45").Value = "BOOL"
$ws.Range("E45").Value = "None"
$ws.Range("F45").Value = "0 - Close, 1 - Open"
$ws.Range("G45").Value = "False"
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1
$ws.Range("L45").Value = "False"

$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "UV_2_OPEN_CLOSE"
$ws.Range("C46").Value = "DB11.DBX12.0"
$ws.Range("D46").Value = "BOOL"
$ws.Range("E46").Value = "None"
$ws.Range("F46").Value = "0 - Close, 1 - Open"
$ws.Range("G46").Value = "False"
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1
$ws.Range("L46").Value = "False"

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "UV_1_BLOCKADE"
$ws.Range("C47").Value = "DB10.DBX12.3"
$ws.Range("D47").Value = "BOOL"
$ws.Range("E47").Value = "None"
$ws.Range("F47").Value = "0 - Close, 1 - Open"
$ws.Range("G47").Value = "False"
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 1
$ws.Range("L47").Value = "False"

$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "UV_2_BLOCKADE"
$ws.Range("C48").Value = "DB11.DBX12.3"
$ws.Range("D48").Value = "BOOL"
$ws.Range("E48").Value = "None"
$ws.Range("F48").Value = "0 - Close, 1 - Open"
$ws.Range("G48").Value = "False"
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 1
$ws.Range("L48").Value = "False"

# --- View state: scroll so row 29 is at the top, select the new rows ---
$ws.Activate()
$excel.Goto($ws.Range("A29"), $true)
$ws.Range("A44:A48").Select()
